# Update the "Förändrad" (Changed) date column (C) for rows 2-15
# from serial date 45174 to 45175 (i.e. 2023-09-05 -> 2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value = 45175
    }
}
